$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 24.21000000000035
$ws.Range("H2").Value = [double]"1.536640864533089e-16"
$ws.Range("K2").Value = 41.42198677273471
$ws.Range("L2").Value = "[36.87763597304462, 45.966337572424806]"
$ws.Range("O2").Value = 1.46544762419704
$ws.Range("P2").Value = "[1.3522370781217328, 1.5786581702723472]"
$ws.Range("S2").Value = 52.61570622894893
$ws.Range("T2").Value = "[49.66277405708723, 55.56863840081063]"
$ws.Range("W2").Value = 18.56342342342369
$ws.Range("X2").Value = 18.12720720720747
$ws.Range("Y2").Value = 18.99963963963991

# Row 3
$ws.Range("E3").Value = 24.97000000000046
$ws.Range("H3").Value = [double]"1.536640864533089e-16"
$ws.Range("K3").Value = 40.25256776716466
$ws.Range("L3").Value = "[34.21914904972262, 46.28598648460671]"
$ws.Range("O3").Value = 2.823974177100735
$ws.Range("P3").Value = "[2.6730267823336575, 2.974921571867812]"
$ws.Range("S3").Value = 52.83386140553301
$ws.Range("T3").Value = "[49.42919607209407, 56.238526738971956]"
$ws.Range("W3").Value = 13.7472472472475
$ws.Range("X3").Value = 13.14736736736761
$ws.Range("Y3").Value = 14.34712712712739
